$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $text) {
    # Writing a literal-looking date/bool string via .Value gets auto-coerced
    # (e.g. "2019-11-04" -> a date serial, "true" -> a boolean) by the
    # engine's type inference, same as typing it into Excel's UI. Routing
    # the text through a formula and then pasting back as a value keeps it
    # a plain shared-string cell (t="s") with no number-format side effects.
    $escaped = $text -replace '"', '""'
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $excel.CutCopyMode = $false
}

# Row 2
Set-TextCell $ws "A2" "ikleiman@stonybrook.edu"
$ws.Range("B2").Value = 1
Set-TextCell $ws "C2" "false"
$ws.Range("D2").Value = 87
Set-TextCell $ws "E2" "2019-11-04"

# Row 3
Set-TextCell $ws "A3" "ikleiman@stonybrook.edu"
$ws.Range("B3").Value = 2
Set-TextCell $ws "C3" "false"
$ws.Range("D3").Value = 106
Set-TextCell $ws "E3" "2019-11-04"

# Row 4
Set-TextCell $ws "A4" "ikleiman@stonybrook.edu"
$ws.Range("B4").Value = 1
Set-TextCell $ws "C4" "false"
$ws.Range("D4").Value = 6
Set-TextCell $ws "E4" "2019-11-04"

# Row 5
Set-TextCell $ws "A5" "ikleiman@stonybrook.edu"
$ws.Range("B5").Value = 2
Set-TextCell $ws "C5" "false"
$ws.Range("D5").Value = 5
Set-TextCell $ws "E5" "2019-11-04"

# Row 6
Set-TextCell $ws "A6" "ikleiman@stonybrook.edu"
$ws.Range("B6").Value = 1
Set-TextCell $ws "C6" "true"
$ws.Range("D6").Value = 103
Set-TextCell $ws "E6" "2019-11-04"

# Row 7
Set-TextCell $ws "A7" "ikleiman@stonybrook.edu"
$ws.Range("B7").Value = 2
Set-TextCell $ws "C7" "true"
$ws.Range("D7").Value = 206
Set-TextCell $ws "E7" "2019-11-04"

# Row 8
Set-TextCell $ws "A8" "ikleiman@stonybrook.edu"
$ws.Range("B8").Value = 1
Set-TextCell $ws "C8" "true"
$ws.Range("D8").Value = 21
Set-TextCell $ws "E8" "2019-11-04"

# Row 9
Set-TextCell $ws "A9" "ikleiman@stonybrook.edu"
$ws.Range("B9").Value = 2
Set-TextCell $ws "C9" "true"
$ws.Range("D9").Value = 9
Set-TextCell $ws "E9" "2019-11-04"

# Row 10
Set-TextCell $ws "A10" "ikleiman@stonybrook.edu"
$ws.Range("B10").Value = 1
Set-TextCell $ws "C10" "true"
$ws.Range("D10").Value = 15
Set-TextCell $ws "E10" "2019-11-04"

# Row 11
Set-TextCell $ws "A11" "ikleiman@stonybrook.edu"
$ws.Range("B11").Value = 2
Set-TextCell $ws "C11" "true"
$ws.Range("D11").Value = 18
Set-TextCell $ws "E11" "2019-11-04"

# Row 12
Set-TextCell $ws "A12" "ikleiman@stonybrook.edu"
$ws.Range("B12").Value = 1
Set-TextCell $ws "C12" "true"
$ws.Range("D12").Value = 15
Set-TextCell $ws "E12" "2019-11-04"

# Row 13
Set-TextCell $ws "A13" "ikleiman@stonybrook.edu"
$ws.Range("B13").Value = 2
Set-TextCell $ws "C13" "true"
$ws.Range("D13").Value = 17
Set-TextCell $ws "E13" "2019-11-04"

# Row 14
Set-TextCell $ws "A14" "ikleiman@stonybrook.edu"
$ws.Range("B14").Value = 1
Set-TextCell $ws "C14" "false"
$ws.Range("D14").Value = 9
Set-TextCell $ws "E14" "2019-11-04"

# Row 15
Set-TextCell $ws "A15" "ikleiman@stonybrook.edu"
$ws.Range("B15").Value = 2
Set-TextCell $ws "C15" "false"
$ws.Range("D15").Value = 27
Set-TextCell $ws "E15" "2019-11-04"

# Row 16
Set-TextCell $ws "A16" "ikleiman@stonybrook.edu"
$ws.Range("B16").Value = 1
Set-TextCell $ws "C16" "true"
$ws.Range("D16").Value = 12
Set-TextCell $ws "E16" "2019-11-04"

# Row 17
Set-TextCell $ws "A17" "ikleiman@stonybrook.edu"
$ws.Range("B17").Value = 2
Set-TextCell $ws "C17" "true"
$ws.Range("D17").Value = 14
Set-TextCell $ws "E17" "2019-11-04"

# Row 18
Set-TextCell $ws "A18" "ikleiman@stonybrook.edu"
$ws.Range("B18").Value = 1
Set-TextCell $ws "C18" "false"
$ws.Range("D18").Value = 4
Set-TextCell $ws "E18" "2019-11-04"

# Row 19
Set-TextCell $ws "A19" "ikleiman@stonybrook.edu"
$ws.Range("B19").Value = 2
Set-TextCell $ws "C19" "false"
$ws.Range("D19").Value = 7
Set-TextCell $ws "E19" "2019-11-04"

# Row 20
Set-TextCell $ws "A20" "ikleiman@stonybrook.edu"
$ws.Range("B20").Value = 1
Set-TextCell $ws "C20" "false"
$ws.Range("D20").Value = 13
Set-TextCell $ws "E20" "2019-11-04"

# Row 21
Set-TextCell $ws "A21" "ikleiman@stonybrook.edu"
$ws.Range("B21").Value = 2
Set-TextCell $ws "C21" "false"
$ws.Range("D21").Value = 15
Set-TextCell $ws "E21" "2019-11-04"

# Row 22
Set-TextCell $ws "A22" "ikleiman@stonybrook.edu"
$ws.Range("B22").Value = 1
Set-TextCell $ws "C22" "false"
$ws.Range("D22").Value = 5
Set-TextCell $ws "E22" "2019-11-04"

# Row 23
Set-TextCell $ws "A23" "ikleiman@stonybrook.edu"
$ws.Range("B23").Value = 2
Set-TextCell $ws "C23" "false"
$ws.Range("D23").Value = 8
Set-TextCell $ws "E23" "2019-11-04"

# Row 24
Set-TextCell $ws "A24" "ikleiman@stonybrook.edu"
$ws.Range("B24").Value = 1
Set-TextCell $ws "C24" "false"
$ws.Range("D24").Value = 6
Set-TextCell $ws "E24" "2019-11-04"

# Row 25
Set-TextCell $ws "A25" "ikleiman@stonybrook.edu"
$ws.Range("B25").Value = 2
Set-TextCell $ws "C25" "true"
$ws.Range("D25").Value = 20
Set-TextCell $ws "E25" "2019-11-04"

# Row 26
Set-TextCell $ws "A26" "chaotsai@stonybrook.edu"
$ws.Range("B26").Value = 1
Set-TextCell $ws "C26" "false"
$ws.Range("D26").Value = 32
Set-TextCell $ws "E26" "2019-12-03"

# Row 27
Set-TextCell $ws "A27" "chaotsai@stonybrook.edu"
$ws.Range("B27").Value = 2
Set-TextCell $ws "C27" "false"
$ws.Range("D27").Value = 42
Set-TextCell $ws "E27" "2019-12-03"

# Row 28
Set-TextCell $ws "A28" "vlgarcia@stonybrook.edu"
$ws.Range("B28").Value = 1
Set-TextCell $ws "C28" "true"
$ws.Range("D28").Value = 73
Set-TextCell $ws "E28" "2019-12-03"

# Row 29
Set-TextCell $ws "A29" "vlgarcia@stonybrook.edu"
$ws.Range("B29").Value = 2
Set-TextCell $ws "C29" "true"
$ws.Range("D29").Value = 80
Set-TextCell $ws "E29" "2019-12-03"
